$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.960.44'
$ws.Range('E2').Value = '  -2.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.412.50'
$ws.Range('E3').Value = '  -2.80%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.10'
$ws.Range('E5').Value = '  -2.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '127.09'
$ws.Range('E6').Value = '  -5.61%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.412.00'
$ws.Range('E8').Value = '  -2.81%  '
$ws.Range('E9').Value = '  -2.46%  '
$ws.Range('E10').Value = '  -1.89%  '
$ws.Range('E11').Value = '  -2.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.381'
$ws.Range('E12').Value = '  -1.42%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.998.56'
$ws.Range('E13').Value = '  -2.69%  '
$ws.Range('E14').Value = '  -0.88%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.416.95'
$ws.Range('E15').Value = '  -2.76%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000174'
$ws.Range('E16').Value = '  -3.97%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.023.76'
$ws.Range('E17').Value = '  -2.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.00'
$ws.Range('E18').Value = '  -3.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.59'
$ws.Range('E19').Value = '  -3.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.71'
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.21'
$ws.Range('E21').Value = '  -3.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '378.50'
$ws.Range('E22').Value = '  -4.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.561'
$ws.Range('E23').Value = '  -2.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.549.78'
$ws.Range('E24').Value = '  -2.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '72.62'
$ws.Range('E25').Value = '  -2.85%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('E27').Value = '  -7.51%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('E29').Value = '  -5.30%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.17'
$ws.Range('E30').Value = '  -4.13%  '
$ws.Range('E31').Value = '  -4.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.42'
$ws.Range('E32').Value = '  -3.53%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.152'
$ws.Range('E33').Value = '  -2.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.441.92'
$ws.Range('E34').Value = '  -2.76%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  -2.08%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.78'
$ws.Range('E38').Value = '  -2.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '164.09'
$ws.Range('E39').Value = '  -1.89%  '
$ws.Range('E40').Value = '  -2.81%  '
$ws.Range('E41').Value = '  -3.23%  '
$ws.Range('E42').Value = '  -3.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.74'
$ws.Range('E44').Value = '  -1.40%  '
$ws.Range('E45').Value = '  -3.32%  '
$ws.Range('E46').Value = '  -5.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '22.98'
$ws.Range('E47').Value = '  -8.63%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.10'
$ws.Range('E48').Value = '  -6.51%  '
$ws.Range('E49').Value = '  -1.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.867'
$ws.Range('E50').Value = '  -3.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.264.57'
$ws.Range('E51').Value = '  -5.15%  '
